$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('SchemaOrganization')
$ws.Range('B6').Value2 = 'http://example.com/organization3:Image0'
$ws.Range('B7').Value2 = 'http://example.com/organization3:Image0'
$ws.Range('B8').Value2 = 'http://example.com/organization3:Image0'
$ws.Range('B9').Value2 = 'http://example.com/organization3:Image0'

$ws = $wb.Worksheets.Item('FoafPerson')
$ws.Range('E2').Value2 = 'http://example.com/person0:Image0'
$ws.Range('E3').Value2 = 'http://example.com/person0:Image0'
$ws.Range('E4').Value2 = 'http://example.com/person0:Image0'
$ws.Range('E5').Value2 = 'http://example.com/person0:Image0'
$ws.Range('E6').Value2 = 'http://example.com/person2:Image0'
$ws.Range('E7').Value2 = 'http://example.com/person2:Image0'
$ws.Range('E8').Value2 = 'http://example.com/person2:Image0'
$ws.Range('E9').Value2 = 'http://example.com/person2:Image0'
$ws.Range('E10').Value2 = 'http://example.com/person4:Image1'
$ws.Range('E11').Value2 = 'http://example.com/person4:Image1'
$ws.Range('E12').Value2 = 'http://example.com/person4:Image1'
$ws.Range('E13').Value2 = 'http://example.com/person4:Image1'

$ws = $wb.Worksheets.Item('SchemaPerson')
$ws.Range('D2').Value2 = 'http://example.com/person1:Image0'
$ws.Range('D3').Value2 = 'http://example.com/person1:Image0'
$ws.Range('D4').Value2 = 'http://example.com/person1:Image0'
$ws.Range('D5').Value2 = 'http://example.com/person1:Image0'

$ws = $wb.Worksheets.Item('SchemaExhibitionEvent')
$ws.Range('B6').Value2 = 'http://example.com/exhibitionEvent1:Image1'
$ws.Range('B7').Value2 = 'http://example.com/exhibitionEvent1:Image1'
$ws.Range('B8').Value2 = 'http://example.com/exhibitionEvent1:Image1'
$ws.Range('B9').Value2 = 'http://example.com/exhibitionEvent1:Image1'
$ws.Range('B14').Value2 = 'http://example.com/exhibitionEvent3:Image1'
$ws.Range('B15').Value2 = 'http://example.com/exhibitionEvent3:Image1'
$ws.Range('B16').Value2 = 'http://example.com/exhibitionEvent3:Image1'
$ws.Range('B17').Value2 = 'http://example.com/exhibitionEvent3:Image1'

$ws = $wb.Worksheets.Item('RdfProperty')
$ws.Range('C10').Value2 = 'dcterms:language:Image1'
$ws.Range('C11').Value2 = 'dcterms:language:Image1'
$ws.Range('C12').Value2 = 'dcterms:language:Image1'
$ws.Range('C13').Value2 = 'dcterms:language:Image1'
$ws.Range('C22').Value2 = 'dcterms:source:Image0'
$ws.Range('C23').Value2 = 'dcterms:source:Image0'
$ws.Range('C24').Value2 = 'dcterms:source:Image0'
$ws.Range('C25').Value2 = 'dcterms:source:Image0'
$ws.Range('C34').Value2 = 'dcterms:title:Image0'
$ws.Range('C35').Value2 = 'dcterms:title:Image0'
$ws.Range('C36').Value2 = 'dcterms:title:Image0'
$ws.Range('C37').Value2 = 'dcterms:title:Image0'
$ws.Range('C38').Value2 = 'dcterms:type:Image1'
$ws.Range('C39').Value2 = 'dcterms:type:Image1'
$ws.Range('C40').Value2 = 'dcterms:type:Image1'
$ws.Range('C41').Value2 = 'dcterms:type:Image1'

$ws = $wb.Worksheets.Item('SchemaProperty')
$ws.Range('C2').Value2 = 'schema:description:Image1'
$ws.Range('C3').Value2 = 'schema:description:Image1'
$ws.Range('C4').Value2 = 'schema:description:Image1'
$ws.Range('C5').Value2 = 'schema:description:Image1'
$ws.Range('C6').Value2 = 'schema:name:Image0'
$ws.Range('C7').Value2 = 'schema:name:Image0'
$ws.Range('C8').Value2 = 'schema:name:Image0'
$ws.Range('C9').Value2 = 'schema:name:Image0'

$ws = $wb.Worksheets.Item('FoafOrganization')
$ws.Range('C2').Value2 = 'http://example.com/organization0:Image0'
$ws.Range('C3').Value2 = 'http://example.com/organization0:Image0'
$ws.Range('C4').Value2 = 'http://example.com/organization0:Image0'
$ws.Range('C5').Value2 = 'http://example.com/organization0:Image0'
$ws.Range('C6').Value2 = 'http://example.com/organization2:Image1'
$ws.Range('C7').Value2 = 'http://example.com/organization2:Image1'
$ws.Range('C8').Value2 = 'http://example.com/organization2:Image1'
$ws.Range('C9').Value2 = 'http://example.com/organization2:Image1'
$ws.Range('C10').Value2 = 'http://example.com/organization4:Image0'
$ws.Range('C11').Value2 = 'http://example.com/organization4:Image0'
$ws.Range('C12').Value2 = 'http://example.com/organization4:Image0'
$ws.Range('C13').Value2 = 'http://example.com/organization4:Image0'

$ws = $wb.Worksheets.Item('SkosConcept')
$ws.Range('B6').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0'
$ws.Range('B7').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0'
$ws.Range('B8').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0'
$ws.Range('B9').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:1:Image0'
$ws.Range('B22').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0'
$ws.Range('B23').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0'
$ws.Range('B24').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0'
$ws.Range('B25').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image0'
$ws.Range('B30').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image0'
$ws.Range('B31').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image0'
$ws.Range('B32').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image0'
$ws.Range('B33').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:7:Image0'
$ws.Range('B34').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0'
$ws.Range('B35').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0'
$ws.Range('B36').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0'
$ws.Range('B37').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0'
$ws.Range('B38').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image1'
$ws.Range('B39').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image1'
$ws.Range('B40').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image1'
$ws.Range('B41').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image1'
$ws.Range('B42').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image0'
$ws.Range('B43').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image0'
$ws.Range('B44').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image0'
$ws.Range('B45').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image0'
$ws.Range('B54').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:13:Image0'
$ws.Range('B55').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:13:Image0'
$ws.Range('B56').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:13:Image0'
$ws.Range('B57').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:13:Image0'
$ws.Range('B58').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0'
$ws.Range('B59').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0'
$ws.Range('B60').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0'
$ws.Range('B61').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:14:Image0'
$ws.Range('B62').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0'
$ws.Range('B63').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0'
$ws.Range('B64').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0'
$ws.Range('B65').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:15:Image0'
$ws.Range('B66').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0'
$ws.Range('B67').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0'
$ws.Range('B68').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0'
$ws.Range('B69').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0'
$ws.Range('B74').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0'
$ws.Range('B75').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0'
$ws.Range('B76').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0'
$ws.Range('B77').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image0'
$ws.Range('B82').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1'
$ws.Range('B83').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1'
$ws.Range('B84').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1'
$ws.Range('B85').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1'
$ws.Range('B94').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1'
$ws.Range('B95').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1'
$ws.Range('B96').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1'
$ws.Range('B97').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:23:Image1'
$ws.Range('B102').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image0'
$ws.Range('B103').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image0'
$ws.Range('B104').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image0'
$ws.Range('B105').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image0'
$ws.Range('B106').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image1'
$ws.Range('B107').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image1'
$ws.Range('B108').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image1'
$ws.Range('B109').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:26:Image1'
$ws.Range('B110').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image1'
$ws.Range('B111').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image1'
$ws.Range('B112').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image1'
$ws.Range('B113').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image1'
$ws.Range('B114').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0'
$ws.Range('B115').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0'
$ws.Range('B116').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0'
$ws.Range('B117').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:28:Image0'
$ws.Range('B122').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1'
$ws.Range('B123').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1'
$ws.Range('B124').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1'
$ws.Range('B125').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1'
$ws.Range('B126').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1'
$ws.Range('B127').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1'
$ws.Range('B128').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1'
$ws.Range('B129').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1'
$ws.Range('B130').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0'
$ws.Range('B131').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0'
$ws.Range('B132').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0'
$ws.Range('B133').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image0'
$ws.Range('B138').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image1'
$ws.Range('B139').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image1'
$ws.Range('B140').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image1'
$ws.Range('B141').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:34:Image1'
$ws.Range('B146').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1'
$ws.Range('B147').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1'
$ws.Range('B148').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1'
$ws.Range('B149').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:36:Image1'
$ws.Range('B162').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0'
$ws.Range('B163').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0'
$ws.Range('B164').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0'
$ws.Range('B165').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:40:Image0'
$ws.Range('B174').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image1'
$ws.Range('B175').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image1'
$ws.Range('B176').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image1'
$ws.Range('B177').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:43:Image1'
$ws.Range('B194').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image1'
$ws.Range('B195').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image1'
$ws.Range('B196').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image1'
$ws.Range('B197').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image1'
$ws.Range('B210').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image0'
$ws.Range('B211').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image0'
$ws.Range('B212').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image0'
$ws.Range('B213').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image0'
$ws.Range('B214').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image1'
$ws.Range('B215').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image1'
$ws.Range('B216').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image1'
$ws.Range('B217').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:53:Image1'
$ws.Range('B218').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image0'
$ws.Range('B219').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image0'
$ws.Range('B220').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image0'
$ws.Range('B221').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:54:Image0'
$ws.Range('B230').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image1'
$ws.Range('B231').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image1'
$ws.Range('B232').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image1'
$ws.Range('B233').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image1'
$ws.Range('B246').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:61:Image1'
$ws.Range('B247').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:61:Image1'
$ws.Range('B248').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:61:Image1'
$ws.Range('B249').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:61:Image1'
$ws.Range('B274').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image1'
$ws.Range('B275').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image1'
$ws.Range('B276').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image1'
$ws.Range('B277').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image1'
$ws.Range('B278').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1'
$ws.Range('B279').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1'
$ws.Range('B280').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1'
$ws.Range('B281').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:69:Image1'
$ws.Range('B286').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0'
$ws.Range('B287').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0'
$ws.Range('B288').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0'
$ws.Range('B289').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0'
$ws.Range('B298').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1'
$ws.Range('B299').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1'
$ws.Range('B300').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1'
$ws.Range('B301').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image1'
$ws.Range('B302').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0'
$ws.Range('B303').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0'
$ws.Range('B304').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0'
$ws.Range('B305').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0'
$ws.Range('B310').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image1'
$ws.Range('B311').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image1'
$ws.Range('B312').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image1'
$ws.Range('B313').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:77:Image1'
$ws.Range('B318').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image1'
$ws.Range('B319').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image1'
$ws.Range('B320').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image1'
$ws.Range('B321').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:79:Image1'

$ws = $wb.Worksheets.Item('SchemaDefinedTerm')
$ws.Range('B6').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image0'
$ws.Range('B7').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image0'
$ws.Range('B8').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image0'
$ws.Range('B9').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:81:Image0'
$ws.Range('B14').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0'
$ws.Range('B15').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0'
$ws.Range('B16').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0'
$ws.Range('B17').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:83:Image0'
$ws.Range('B18').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image1'
$ws.Range('B19').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image1'
$ws.Range('B20').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image1'
$ws.Range('B21').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:84:Image1'
$ws.Range('B30').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image0'
$ws.Range('B31').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image0'
$ws.Range('B32').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image0'
$ws.Range('B33').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:87:Image0'
$ws.Range('B34').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0'
$ws.Range('B35').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0'
$ws.Range('B36').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0'
$ws.Range('B37').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:88:Image0'
$ws.Range('B38').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image1'
$ws.Range('B39').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image1'
$ws.Range('B40').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image1'
$ws.Range('B41').Value2 = 'urn:paradicms_etl:pipeline:synthetic_data:concept:89:Image1'

$ws = $wb.Worksheets.Item('SchemaCreativeWork')
$ws.Range('E2').Value2 = 'https://images.metmuseum.org/CRDImages/ep/original/LC-EP_1993_132_suppl_CH-001.jpg'
$ws.Range('E3').Value2 = 'https://images.metmuseum.org/CRDImages/ep/original/LC-EP_1993_132_suppl_CH-001.jpg'
$ws.Range('E4').Value2 = 'https://images.metmuseum.org/CRDImages/ep/original/LC-EP_1993_132_suppl_CH-001.jpg'
$ws.Range('E5').Value2 = 'https://images.metmuseum.org/CRDImages/ep/original/LC-EP_1993_132_suppl_CH-001.jpg'

$ws = $wb.Worksheets.Item('RightsStatementsDotOrgRightsStatement')
$ws.Range('A3').Value2 = 'https://rightsstatements.org/vocab/NoC-US/1.0/'
$ws.Range('E3').Value2 = 'You may find additional information about the copyright status of the Item on the website of the organization that has made the Item available.'
$ws.Range('A4').Value2 = 'http://rightsstatements.org/vocab/InC-EDU/1.0/'
$ws.Range('B4').Value2 = 'This Item is protected by copyright and/or related rights.

  You are free to use this Item in any way that is permitted by the copyright and related rights legislation that applies to your use. In addition, no permission is required from the rights-holder(s) for educational uses.

  For other uses, you need to obtain permission from the rights-holder(s).'
$ws.Range('C4').Value2 = 'This Rights Statement indicates that the Item labeled with this Rights Statement is in copyright but that educational use is allowed without the need to obtain additional permission.'
$ws.Range('D4').Value2 = 'InC-EDU'
$ws.Range('F4').Value2 = 'In Copyright - Educational Use Permitted'
$ws.Range('G4').Value2 = 'This Rights Statement can be used only for copyrighted Items for which the organization making the Item available is the rights-holder or has been explicitly authorized by the rights-holder(s) to allow third parties to use their Work(s) for educational purposes without first obtaining permission.'
$ws.Range('A5').Value2 = 'http://rightsstatements.org/vocab/NoC-US/1.0/'
$ws.Range('B5').Value2 = 'The organization that has made the Item available believes that the Item is in the Public Domain under the laws of the United States, but a determination was not made as to its copyright status under the copyright laws of other countries. The Item may not be in the Public Domain under the laws of other countries.

  Please refer to the organization that has made the Item available for more information.'
$ws.Range('C5').Value2 = 'This Rights Statement indicates that the Item is in the Public Domain under the laws of the United States, but that a determination was not made as to its copyright status under the copyright laws of other countries.'
$ws.Range('D5').Value2 = 'NoC-US'
$ws.Range('F5').Value2 = 'No Copyright - United States'
$ws.Range('G5').Value2 = 'This Rights Statement should be used for Items for which the organization that intends to make the Item available has determined are free of copyright under the laws of the United States. This Rights Statement should not be used for Orphan Works (which are assumed to be in-copyright) or for Items where the organization that intends to make the Item available has not undertaken an effort to ascertain the copyright status of the underlying Work.'
